$d = $word.ActiveDocument

$find = " the Delorean dashboard in " + [char]0x201C + "Back to the Future" + [char]0x201D + "."
$replace = " the Delorean dashboard in " + [char]0x201C + "Back to the Future" + [char]0x201D + "."

$d.Content.Find.Execute($find, $true, $false, $false, $false, $false,
                         $true, 1, $false, $replace, 2)
